$d = $word.ActiveDocument

# --- First paragraph: "**ID__AFFARS_pgi_5317_topic_11__ID**" placeholder ---
$p1 = $d.Paragraphs(1)

$oldId = "**ID__AFFARS_pgi_5317_topic_11__ID**"
$newId = "**ID__AFFARS_AFMC_PGI_5317_7506_90__ID**"

# Locate the ID placeholder run so we can find the lone-space run that
# immediately trails it (same run formatting, separate <w:r>) and drop it.
$idRange = $p1.Range.Duplicate
$idRange.Find.Execute($oldId, $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

$trailingSpace = $d.Range($idRange.End, $idRange.End + 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# Update the placeholder text itself.
$p1.Range.Find.Execute($oldId, $false, $false, $false, $false, $false,
    $true, 1, $false, $newId, 2) | Out-Null

# Add a paragraph border (5-twip spacing on all sides) and widen the
# left indent from 120 to 225 twips (6pt -> 11.25pt).
$pFmt = $p1.Format
$pFmt.Borders.DistanceFromTop = 5
$pFmt.Borders.DistanceFromLeft = 5
$pFmt.Borders.DistanceFromBottom = 5
$pFmt.Borders.DistanceFromRight = 5
$pFmt.LeftIndent = 11.25
